$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Table_1")

# --- Table_1: drop the stray empty placeholder cells (B2, A3, B37) ---
$ws1.Range("B2").Value  = $null
$ws1.Range("A3").Value  = $null
$ws1.Range("B37").Value = $null

# --- Add the new Table_2 worksheet right after Table_1 ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Table_2"

# Header row (row 1) text
$ws2.Range("A1").Value = "Əmsal"
$ws2.Range("B1").Value = "Norma (Sistem əhəmiyyətli)"
$ws2.Range("C1").Value = "Norma (Banklar istisna)"
$ws2.Range("D1").Value = "Fakt"

# Reuse Table_1's header formatting (bold, thin border, centered/top-aligned)
# so the new header row shares the same style index instead of creating a new one.
$ws1.Range("A1:B1").Copy()
$ws2.Range("A1:D1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Data rows 2-4: force everything (including "%"-looking text) to be kept as
# literal text rather than auto-converted to a percentage number by pre-marking
# the range as Text, assigning the values, then clearing the format again so no
# residual style sticks to the cells (matches the source which has plain,
# unstyled text cells here).
$dataRng = $ws2.Range("B2:D4")
$dataRng.NumberFormat = "@"

$ws2.Range("A2").Value = "9.  I dərəcəli  kapitalın  adekvatlıq əmsalı"
$ws2.Range("B2").Value = "6.0%"
$ws2.Range("C2").Value = "5.0%"
$ws2.Range("D2").Value = "8.8%"

$ws2.Range("A3").Value = "10. məcmu kapitalın  adekvatlıq  əmsalı"
$ws2.Range("B3").Value = "12.0%"
$ws2.Range("C3").Value = "10.0%"
$ws2.Range("D3").Value = "14.4%"

$ws2.Range("A4").Value = "11. Leverec əmsalı"
$ws2.Range("B4").Value = "minimum 5%"
$ws2.Range("C4").Value = "minimum 4%"
$ws2.Range("D4").Value = "5.25%"

$dataRng.ClearFormats()

# Restore Table_1 as the active sheet (matches the original workbook's state).
$ws1.Activate()
